$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row for the "2_skin" alien feature, right after the
# existing "2_body" row (original row 39, before any other shifts).
$ws.Rows("39:39").Insert()
$ws.Range("A39").Value = "2_skin"
$ws.Range("B39").Value = "Data/Images/Features_Placed/2_skin.png"

# Insert a new row for the "1_skin" alien feature, right after the
# existing "1_body" row (row 2).
$ws.Rows("3:3").Insert()
$ws.Range("A3").Value = "1_skin"
$ws.Range("B3").Value = "Data/Images/Features_Placed/1_skin.png"

# Move the active selection to B4 (matches the saved view state).
$ws.Range("B4").Select()
